$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataCellStyle($rng) {
  $rng.Borders.LineStyle = 1
  $rng.Borders.Color = 0
  $rng.Font.Name = "TimesNewRoman"
  $rng.Font.Bold = $true
  $rng.Font.Size = 12
  $rng.Font.Color = 0
  $rng.HorizontalAlignment = -4108
  $rng.VerticalAlignment = -4108
}

# ---- Row 5: update existing cells (C5, W5) ----
$ws.Range("C5").Value = "P1"
$ws.Range("W5").FormulaLocal = "'3"

# ---- Row 6: update existing cells ----
$ws.Range("A6").Value = "Natasha"
$ws.Range("C6").Value = "A1"
$ws.Range("E6").Value = "F"
$ws.Range("G6").Value = "I"
$ws.Range("I6").FormulaLocal = "'1"
$ws.Range("K6").Value = 230
$ws.Range("M6").Value = "E"
$ws.Range("O6").FormulaLocal = "'4"
$ws.Range("Q6").Value = 0
$ws.Range("S6").Value = 100
$ws.Range("U6").Value = "C"
$ws.Range("W6").FormulaLocal = "'3"
$ws.Range("Y6").Value = "R"
$ws.Range("AA6").Value = 16

# ---- Row 7: update existing cells ----
$ws.Range("A7").Value = "Popka"
$ws.Range("C7").Value = "A3"
$ws.Range("E7").Value = "F"
$ws.Range("G7").Value = "I"
$ws.Range("I7").FormulaLocal = "'1"
$ws.Range("K7").Value = 230
$ws.Range("M7").Value = "E"
$ws.Range("O7").FormulaLocal = "'4"
$ws.Range("Q7").Value = 43
$ws.Range("S7").Value = 100
$ws.Range("U7").Value = "C"
$ws.Range("W7").FormulaLocal = "'3"
$ws.Range("Y7").Value = "R"
$ws.Range("AA7").Value = 16

# ---- Row 8: update existing cells ----
$ws.Range("A8").Value = "Sisechki"
$ws.Range("C8").Value = "A4"
$ws.Range("E8").Value = "F"
$ws.Range("G8").Value = "III"
$ws.Range("I8").FormulaLocal = "'3"
$ws.Range("K8").Value = 528
$ws.Range("M8").Value = "E"
$ws.Range("O8").FormulaLocal = "'5"
$ws.Range("Q8").Value = 56
$ws.Range("S8").Value = 304
$ws.Range("U8").Value = "D"
$ws.Range("W8").Value = "N"
$ws.Range("Y8").Value = "I"
$ws.Range("AA8").Value = 20

# ---- Row 9: update existing cells ----
$ws.Range("A9").Value = "Ssadasd"
$ws.Range("C9").Value = "A3"
$ws.Range("E9").Value = "H"
$ws.Range("G9").Value = "IV"
$ws.Range("I9").FormulaLocal = "'3"
$ws.Range("K9").Value = 1216
$ws.Range("M9").Value = "H"
$ws.Range("O9").FormulaLocal = "'5"
$ws.Range("Q9").Value = 87
$ws.Range("S9").Value = 1056
$ws.Range("U9").Value = "H"
$ws.Range("W9").FormulaLocal = "'4"
$ws.Range("Y9").Value = "S"
$ws.Range("AA9").Value = 28

# ---- Row 10: update existing cells ----
$ws.Range("A10").Value = "Safsdf"
$ws.Range("C10").Value = "A3"
$ws.Range("E10").Value = "F"
$ws.Range("G10").Value = "I"
$ws.Range("I10").FormulaLocal = "'1"
$ws.Range("K10").Value = 230
$ws.Range("M10").Value = "E"
$ws.Range("O10").FormulaLocal = "'4"
$ws.Range("Q10").Value = 43
$ws.Range("S10").Value = 100
$ws.Range("U10").Value = "C"
$ws.Range("W10").FormulaLocal = "'3"
$ws.Range("Y10").Value = "R"
$ws.Range("AA10").Value = 16

# ---- Row 11: new row ----
$ws.Range("A11:B11").Merge()
$ws.Range("C11:D11").Merge()
$ws.Range("E11:F11").Merge()
$ws.Range("G11:H11").Merge()
$ws.Range("I11:J11").Merge()
$ws.Range("K11:L11").Merge()
$ws.Range("M11:N11").Merge()
$ws.Range("O11:P11").Merge()
$ws.Range("Q11:R11").Merge()
$ws.Range("S11:T11").Merge()
$ws.Range("U11:V11").Merge()
$ws.Range("W11:X11").Merge()
$ws.Range("Y11:Z11").Merge()
Set-DataCellStyle($ws.Range("A11:B11"))
Set-DataCellStyle($ws.Range("C11:D11"))
Set-DataCellStyle($ws.Range("E11:F11"))
Set-DataCellStyle($ws.Range("G11:H11"))
Set-DataCellStyle($ws.Range("I11:J11"))
Set-DataCellStyle($ws.Range("K11:L11"))
Set-DataCellStyle($ws.Range("M11:N11"))
Set-DataCellStyle($ws.Range("O11:P11"))
Set-DataCellStyle($ws.Range("Q11:R11"))
Set-DataCellStyle($ws.Range("S11:T11"))
Set-DataCellStyle($ws.Range("U11:V11"))
Set-DataCellStyle($ws.Range("W11:X11"))
Set-DataCellStyle($ws.Range("Y11:Z11"))
Set-DataCellStyle($ws.Range("AA11"))
$ws.Range("A11").Value = "Zhopa"
$ws.Range("C11").Value = "A3"
$ws.Range("E11").Value = "D"
$ws.Range("G11").Value = "I"
$ws.Range("I11").FormulaLocal = "'3"
$ws.Range("K11").Value = 175
$ws.Range("M11").Value = "B"
$ws.Range("O11").FormulaLocal = "'2"
$ws.Range("Q11").Value = 16
$ws.Range("S11").Value = 29
$ws.Range("U11").Value = "A"
$ws.Range("W11").FormulaLocal = "'4"
$ws.Range("Y11").Value = "P"
$ws.Range("AA11").Value = 14

# ---- Row 12: new row ----
$ws.Range("A12:B12").Merge()
$ws.Range("C12:D12").Merge()
$ws.Range("E12:F12").Merge()
$ws.Range("G12:H12").Merge()
$ws.Range("I12:J12").Merge()
$ws.Range("K12:L12").Merge()
$ws.Range("M12:N12").Merge()
$ws.Range("O12:P12").Merge()
$ws.Range("Q12:R12").Merge()
$ws.Range("S12:T12").Merge()
$ws.Range("U12:V12").Merge()
$ws.Range("W12:X12").Merge()
$ws.Range("Y12:Z12").Merge()
Set-DataCellStyle($ws.Range("A12:B12"))
Set-DataCellStyle($ws.Range("C12:D12"))
Set-DataCellStyle($ws.Range("E12:F12"))
Set-DataCellStyle($ws.Range("G12:H12"))
Set-DataCellStyle($ws.Range("I12:J12"))
Set-DataCellStyle($ws.Range("K12:L12"))
Set-DataCellStyle($ws.Range("M12:N12"))
Set-DataCellStyle($ws.Range("O12:P12"))
Set-DataCellStyle($ws.Range("Q12:R12"))
Set-DataCellStyle($ws.Range("S12:T12"))
Set-DataCellStyle($ws.Range("U12:V12"))
Set-DataCellStyle($ws.Range("W12:X12"))
Set-DataCellStyle($ws.Range("Y12:Z12"))
Set-DataCellStyle($ws.Range("AA12"))
$ws.Range("A12").Value = "Pisya"
$ws.Range("C12").Value = "A3"
$ws.Range("E12").Value = "C"
$ws.Range("G12").Value = "III"
$ws.Range("I12").FormulaLocal = "'2"
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = "B"
$ws.Range("O12").FormulaLocal = "'4"
$ws.Range("Q12").Value = 28
$ws.Range("S12").Value = 57
$ws.Range("U12").Value = "B"
$ws.Range("W12").FormulaLocal = "'3"
$ws.Range("Y12").Value = "C"
$ws.Range("AA12").Value = 14

